# edit.ps1 -- apply README.docx changes via Word COM-interop
$d = $word.ActiveDocument

function Replace-Exact($oldText, $newText) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $found = $rng.Find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $found) {
        throw "Find.Execute could not find: $oldText"
    }
}

# ---------------------------------------------------------------------
# 1) "1. .cfg files" / "2. .gpf files" / "3. Finding MESes from GPFs"
#    -- collapse proofErr-split runs into single clean runs (no visible
#       text change, just tidy-up that happened alongside later edits)
# ---------------------------------------------------------------------
Replace-Exact "1. .cfg files" "1. .cfg files"
Replace-Exact "2. .gpf files" "2. .gpf files"
Replace-Exact "3. Finding MESes from GPFs" "3. Finding MESes from GPFs"

# "1. .cfg files / Configuration Files"
Replace-Exact "1. .cfg files / Configuration Files" "1. .cfg files / Configuration Files"

# ---------------------------------------------------------------------
# 2) Paragraph about hydrophobic center distances / allowance constant
#    contains straight quotes -- use direct Range.Text assignment so
#    Word's smart-quote autocorrect doesn't curl them.
# ---------------------------------------------------------------------
$oldText = 'A .cfg file lets you control the hydrophobic center distances for each hydrophobic residue, as well as the "allowance constant," or AC. The AC is the fraction of points that will be allowed to fall outside of the sphere during Minimal Enclosing Sphere (MES) generation, represented in decimal. That is to say, if you want to make sure that 90% of points fall inside the MES, you would set the AC to .1, so 10% of points will be allowed to fall outside of the MES. If you want 80% inside, you would set the AC to .2.'
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd("`r`a") -eq $oldText) {
        $r = $p.Range
        $r2 = $d.Range($r.Start, $r.End - 1)
        $r2.Text = $oldText
        $found = $true
        break
    }
}
if (-not $found) { throw "Could not locate hydrophobic-center paragraph" }

# ---------------------------------------------------------------------
# 3) "A .cfg file also allows you ... specific identity of each" --
#    collapse the proofErr-split runs (cfg / cender / gpf) into one run
# ---------------------------------------------------------------------
Replace-Exact "A .cfg file also allows you to control the assumed hydrophobic cender distances along the alpha-beta carbon vector. The MES generation script takes in a .gpf file as input, which only has data about the locations of the alpha-carbon backbone, as well as the beta carbons of hydrophobic residues and the specific identity of each" "A .cfg file also allows you to control the assumed hydrophobic cender distances along the alpha-beta carbon vector. The MES generation script takes in a .gpf file as input, which only has data about the locations of the alpha-carbon backbone, as well as the beta carbons of hydrophobic residues and the specific identity of each"

# ---------------------------------------------------------------------
# 4) New paragraph ".cfg files are stored in the config folder." right
#    after "... distance is determined by the specific identity of the
#    residue."
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd("`r`a").EndsWith("distance is determined by the specific identity of the residue.")) {
        $newPara = $p.Range.InsertParagraphAfter()
        $d.Paragraphs.Item($i + 1).Range.Text = ".cfg files are stored in the config folder."
        break
    }
}

# ---------------------------------------------------------------------
# 5) "A .cfg file is very simple to make. ... Each line has" -- collapse
# ---------------------------------------------------------------------
Replace-Exact "A .cfg file is very simple to make. It is a series of lines of text. Each line has" "A .cfg file is very simple to make. It is a series of lines of text. Each line has"

# ---------------------------------------------------------------------
# 6) "allowance constant is set to .2. To see a full example, check the
#    default.cfg" -- collapse proofErr-split run
# ---------------------------------------------------------------------
Replace-Exact "allowance constant is set to .2. To see a full example, check the default.cfg" "allowance constant is set to .2. To see a full example, check the default.cfg"

# ---------------------------------------------------------------------
# 7) "2. .gpf files / Geometric Protein Files" -- collapse
# ---------------------------------------------------------------------
Replace-Exact "2. .gpf files / Geometric Protein Files" "2. .gpf files / Geometric Protein Files"

# ---------------------------------------------------------------------
# 8) "1. What and why?" -> "2.1. What and why?"
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd("`r`a") -eq "1. What and why?") {
        $r = $p.Range
        $r2 = $d.Range($r.Start, $r.End - 1)
        $r2.Text = "2.1. What and why?"
        break
    }
}

# ---------------------------------------------------------------------
# 9) GPF paragraph: append "folder in the script directory." after the
#    existing "... placed in the gpf_files " text.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("A GPF consists of two components.")) {
        $endRange = $d.Range($p.Range.End - 1, $p.Range.End - 1)
        $endRange.InsertAfter("folder in the script directory.")
        $gpfParaIndex = $i
        break
    }
}

# ---------------------------------------------------------------------
# 10) New paragraphs appended after the GPF paragraph:
#     "2.2. How?"
#     "At the moment, .gpf files can only be generated by running the
#      GeoProtCreate.py script. Simply give it a PDB protein name and
#      it'll do the rest for you."
#     "3. Finding MESes from GPFs"
#     "Run the MESFromGPF.py script and follow the instructions."
#     <empty paragraph>
# ---------------------------------------------------------------------
$gpfPara = $d.Paragraphs.Item($gpfParaIndex)
$cursor = $gpfPara.Range.InsertParagraphAfter()
$d.Paragraphs.Item($gpfParaIndex + 1).Range.Text = "2.2. How?"

$p = $d.Paragraphs.Item($gpfParaIndex + 1)
$p.Range.InsertParagraphAfter()
$d.Paragraphs.Item($gpfParaIndex + 2).Range.Text = "At the moment, .gpf files can only be generated by running the GeoProtCreate.py script. Simply give it a PDB protein name and it" + [char]0x2019 + "ll do the rest for you."

$p = $d.Paragraphs.Item($gpfParaIndex + 2)
$p.Range.InsertParagraphAfter()
$d.Paragraphs.Item($gpfParaIndex + 3).Range.Text = "3. Finding MESes from GPFs"

$p = $d.Paragraphs.Item($gpfParaIndex + 3)
$p.Range.InsertParagraphAfter()
$d.Paragraphs.Item($gpfParaIndex + 4).Range.Text = "Run the MESFromGPF.py script and follow the instructions."

$p = $d.Paragraphs.Item($gpfParaIndex + 4)
$p.Range.InsertParagraphAfter()

Write-Output "Done"
